$wb = $excel.ActiveWorkbook

# Update "zh-cn" sheet: row 2 (0db6ae1f... file) handoff/handback datetimes were refreshed
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-28 11:02:23"
$wsZhCn.Range("G2").Value = "2016-01-28 11:03:10"

# Update "de-de" sheet: row 2 (0db6ae1f... file) handoff/handback datetimes were refreshed
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-28 11:02:40"
$wsDeDe.Range("G2").Value = "2016-01-28 11:03:32"
